$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{}
$data["B"] = @(12.77268770170907,12.28844379339649,11.98255873707406,11.85594305670654,11.83480578874713,11.98085884854871,12.60760238019068,13.76153456963656,14.5551893360526,14.90306539094039,15.03280613820753,15.00495420349516,14.91377969376548,14.8576704794704,14.53218113448668,14.32905761750773,14.21099256471294,14.17080912926519,14.35080891992601,14.94061465779303,15.31443673450954,15.11601637394322,14.34097915596245,13.45833705309724)
$data["C"] = @(8.008463659254859,7.586567973294963,7.315542951004718,7.202197511780086,7.183205080059635,7.314025916031093,7.865541220452677,8.848056819868191,9.505375168991534,9.78977009321599,9.895321555292609,9.872685008286552,9.798497019048117,9.752774640461611,9.486491436027194,9.3193589428026,9.221854218464074,9.188606075284873,9.337293057572857,9.820346261101747,10.1235469628845,9.962878185970688,9.329189458028726,8.593344972594984)
$data["D"] = @(9.467861914974691,9.457719728724125,9.452769911162383,9.45107592091407,9.450814200735421,9.452745754855696,9.464100896948509,9.496417740375165,9.526162329209077,9.540967197720265,9.546753788058526,9.545499570534805,9.541439667388584,9.538976252822264,9.525220188020361,9.517105567755124,9.512558341148779,9.511039438762484,9.517956974599542,9.542627289559643,9.559800295244035,9.550539697787192,9.517571685873163,9.486611571696962)
$data["E"] = @(13.62946272699483,13.64549136321961,13.65750522157899,13.66294735658101,13.66388402241109,13.65757640365436,13.63453869026721,13.60658677707743,13.5965324558587,13.59422815910757,13.59368114416941,13.59378448745635,13.59417663581345,13.59445921187717,13.59672865137266,13.59870163475667,13.60005011004723,13.60054339960226,13.59846949871771,13.59405262322757,13.59306317044509,13.59341795348192,13.59857378021666,13.61230591894635)
$data["F"] = @(30.704331109877,30.72502909154314,30.74583083609299,30.75634032176948,30.75820808461284,30.74596434485029,30.70978688691693,30.70312280218187,30.73743428077262,30.76153681830931,30.77188140559849,30.76959944239131,30.76236350079856,30.75808969033578,30.73602973628601,30.72466979058258,30.71893566975839,30.7171316646506,30.72579632251312,30.76445586801838,30.79681578889582,30.77889725296061,30.72528453494909,30.6980371149369)
$data["I"] = @(20.51028001455472,20.60406842943532,20.66612236587302,20.69253136154953,20.69698423763949,20.66647398832959,20.54168990985177,20.33251997102021,20.20061897495418,20.14536765739793,20.1251308556785,20.12945867793556,20.14368900740047,20.15249486367334,20.20432559348486,20.23734068947757,20.25677714683831,20.26343470332298,20.2337798905919,20.13949058591613,20.08186486465045,20.11225417656551,20.23538830864397,20.38529031445392)
$data["J"] = @(9.899118514847585,9.92371257768349,9.939934116059327,9.946826644733449,9.947988192524104,9.940025928398628,9.907366139145191,9.852198443268485,9.817059750133954,9.80224125876796,9.79679729303005,9.79796230363789,9.801790025303534,9.804156418692024,9.818051621804429,9.826874422917706,9.83205885399769,9.833833076862513,9.825923859605547,9.800661187281037,9.785126638568691,9.793328485767052,9.826353260123055,9.866174316482107)
$data["M"] = @(16.36564763522155,16.22393802648412,16.13878318282086,16.10457990110749,16.0989314525612,16.13831984809511,16.31642018386841,16.67897920019148,16.95152496174483,17.07643441589403,17.12383231464361,17.1136205671083,17.0803321225941,17.05995361050252,16.94337741684663,16.87207354727391,16.83115084836158,16.81731157641157,16.87965497377822,17.09010738760552,17.22820141960694,17.15445941121228,16.87622718743029,16.57967132168922)
$data["N"] = @(17.44350750119754,17.49257171515311,17.52444956731046,17.53788156545414,17.54013863172739,17.52462892705671,17.46006170892699,17.3473079766124,17.27286232805365,17.24080618717757,17.22892673254165,17.23147365481156,17.23982366026928,17.24497205513272,17.27499362467987,17.29387385010865,17.3049036211629,17.30866739418883,17.29184639091331,17.23736402371481,17.20326896075478,17.22132798713882,17.29276245916976,17.37633269123383)
$data["O"] = @(22.8532742165727,22.89958564744176,22.93369522678467,22.94901834871429,22.95164858460267,22.9338961224392,22.86806266419976,22.7841309005965,22.7501751735901,22.74077200883927,22.73808168012685,22.73862236064705,22.7405332227313,22.74181707109294,22.75091143466527,22.75803957671052,22.76270828476727,22.76438666771556,22.75722189437795,22.73994832406407,22.73373293858841,22.73658562189202,22.75758979113068,22.8019813695105)

$startRow = 2
foreach ($col in $data.Keys) {
    $colIndex = $ws.Range($col + "1").Column
    $vals = $data[$col]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($startRow + $i, $colIndex).Value = $vals[$i]
    }
}
